# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps for the
# a73cac21-... handback entry on both the zh-cn and de-de sheets.
# Rows 3 and 5 on each sheet shared the same datetime text, so both
# rows are updated to keep them sharing the (new) value, matching
# the original shared-string reuse.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-20 02:16:54"
$wsZhCn.Range("H3").Value = "2016-03-20 02:17:13"
$wsZhCn.Range("E5").Value = "2016-03-20 02:16:54"
$wsZhCn.Range("H5").Value = "2016-03-20 02:17:13"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-20 02:16:57"
$wsDeDe.Range("H3").Value = "2016-03-20 02:17:19"
$wsDeDe.Range("E5").Value = "2016-03-20 02:16:57"
$wsDeDe.Range("H5").Value = "2016-03-20 02:17:19"
